$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 19: 12/03/2025, 1:00 hours, task note ---
# Copy formats from row 18 (same date/hours column styles) so no new
# number-format/style entries are introduced, then overwrite the values.
# A and C are copied individually (not A:C) so an empty B cell isn't
# materialised, matching the target sheet which only has A/C/E on row 19.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("A19").Value = [DateTime]"2025-12-03"
$ws.Range("C19").Value = 0.041666666666666664
$ws.Range("E19").Value = "created this, talked o you, and more"

# --- Row 21: note only in column A ---
# Written before row 20's text so the shared-string table fills in the same
# order as the target file: 23=created this..., 24=pips..., 25=meeting with prof.
$ws.Range("A21").Value = "pips: gymnasium, ale_py, torch,"

# --- Row 20: 12/09/2025, 0:20 hours, task note ---
$ws.Range("A18").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("C18").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("A20").Value = [DateTime]"2025-12-09"
$ws.Range("C20").Value = 0.013888888888888888
$ws.Range("E20").Value = "meeting with prof."

$excel.CutCopyMode = 0

# --- Update the view / selection to match the saved state ---
$ws.Range("E20").Select()

$wb.Save()
